$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D19","D20","D22","D23","D24","D25","D26","D28","D31","D32","D33","D34","D35","D36","D37","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.733.51"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.947.90"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "247.44"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.4842"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").Value = "0.2951"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").Value = "0.06821"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "112.39"
$ws.Range("E10").Value = "  +5.06%  "
$ws.Range("D11").Value = "19.52"
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").Value = "1.940.12"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "5.543"
$ws.Range("E13").Value = "  +5.44%  "
$ws.Range("D14").Value = "0.07649"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "0.6918"
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("D16").Value = "295.97"
$ws.Range("E16").Value = "  +9.75%  "
$ws.Range("D17").Value = "30.723.84"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "13.37"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("D19").Value = "5.701"
$ws.Range("D20").Value = "0.000007701"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").Value = "2.198.35"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "6.558"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").Value = "9.774"
$ws.Range("E25").Value = "  +4.13%  "
$ws.Range("D26").Value = "167.91"
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("D28").Value = "2.179"
$ws.Range("E28").Value = "  +3.91%  "
$ws.Range("E29").Value = "  +4.10%  "
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").Value = "4.781"
$ws.Range("E31").Value = "  +18.58%  "
$ws.Range("D32").Value = "4.433"
$ws.Range("E32").Value = "  +7.72%  "
$ws.Range("D33").Value = "0.05074"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").Value = "0.7786"
$ws.Range("E34").Value = "  +7.22%  "
$ws.Range("D35").Value = "1.162"
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "0.02071"
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("D37").Value = "2.736"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").Value = "2.044"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").Value = "110.98"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "0.4454"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("D42").Value = "0.8743"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "5.955"
$ws.Range("E43").Value = "  +1.71%  "
$ws.Range("D44").Value = "70.95"
$ws.Range("E44").Value = "  +5.00%  "
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "7.389"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "9.462"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("D48").Value = "48.93"
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("D49").Value = "0.1251"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "35.52"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("D51").Value = "0.2550"
$ws.Range("E51").Value = "  +2.96%  "
